$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("mean")
$ws.Range("B2").Value = 70.42100000000001
$ws.Range("C2").Value = 19.286
$ws.Range("D2").Value = 0.031
$ws.Range("E2").Value = 0.126
$ws.Range("F2").Value = 1.232
$ws.Range("G2").Value = 6.726
$ws.Range("H2").Value = 97.82299999999999
$ws.Range("B3").Value = 70.464
$ws.Range("C3").Value = 19.3
$ws.Range("D3").Value = 0.031
$ws.Range("E3").Value = 0.126
$ws.Range("F3").Value = 1.266
$ws.Range("G3").Value = 6.908
$ws.Range("H3").Value = 98.096
$ws.Range("B4").Value = 70.806
$ws.Range("C4").Value = 19.408
$ws.Range("F4").Value = 1.533
$ws.Range("G4").Value = 8.369
$ws.Range("H4").Value = 100.273

$ws = $wb.Worksheets.Item("stdev")
$ws.Range("B2").Value = 0.489
$ws.Range("C2").Value = 0.12
$ws.Range("E2").Value = 0.033
$ws.Range("F2").Value = 0.056
$ws.Range("G2").Value = 0.304
$ws.Range("H2").Value = 0.722
$ws.Range("B3").Value = 0.487
$ws.Range("C3").Value = 0.12
$ws.Range("E3").Value = 0.033
$ws.Range("F3").Value = 0.054
$ws.Range("G3").Value = 0.294
$ws.Range("H3").Value = 0.708
$ws.Range("B4").Value = 0.503
$ws.Range("C4").Value = 0.122
$ws.Range("D4").Value = 0.012
$ws.Range("F4").Value = 0.065
$ws.Range("G4").Value = 0.354
$ws.Range("H4").Value = 0.787

$ws = $wb.Worksheets.Item("summary")
$ws.Range("B2").Value = 19.286
$ws.Range("C2").Value = 0.12
$ws.Range("D2").Value = 1.232
$ws.Range("E2").Value = 0.056
$ws.Range("F2").Value = 0.031
$ws.Range("H2").Value = 6.726
$ws.Range("I2").Value = 0.304
$ws.Range("J2").Value = 0.126
$ws.Range("K2").Value = 0.033
$ws.Range("L2").Value = 70.42100000000001
$ws.Range("M2").Value = 0.489
$ws.Range("N2").Value = 97.82299999999999
$ws.Range("O2").Value = 0.722
$ws.Range("B3").Value = 19.3
$ws.Range("C3").Value = 0.12
$ws.Range("D3").Value = 1.266
$ws.Range("E3").Value = 0.054
$ws.Range("F3").Value = 0.031
$ws.Range("H3").Value = 6.908
$ws.Range("I3").Value = 0.294
$ws.Range("J3").Value = 0.126
$ws.Range("K3").Value = 0.033
$ws.Range("L3").Value = 70.464
$ws.Range("M3").Value = 0.487
$ws.Range("N3").Value = 98.096
$ws.Range("O3").Value = 0.708
$ws.Range("B4").Value = 19.408
$ws.Range("C4").Value = 0.122
$ws.Range("D4").Value = 1.533
$ws.Range("E4").Value = 0.065
$ws.Range("G4").Value = 0.012
$ws.Range("H4").Value = 8.369
$ws.Range("I4").Value = 0.354
$ws.Range("L4").Value = 70.806
$ws.Range("M4").Value = 0.503
$ws.Range("N4").Value = 100.273
$ws.Range("O4").Value = 0.787
